$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 16917.715
$ws.Range("J95").Value = 16917.715
$ws.Range("L95").Value = 16917.715
$ws.Range("N95").Value = -22409.715

$ws.Range("H132").Value = 2315.92
$ws.Range("I132").Value = 1222.6818
$ws.Range("J132").Value = 10333
$ws.Range("K132").Value = 3668.0454
$ws.Range("L132").Value = 30999
$ws.Range("M132").Value = -1138.0454
$ws.Range("N132").Value = -36059

$ws.Range("H137").Value = 2334.375
$ws.Range("I137").Value = 2095.8
$ws.Range("K137").Value = 6287.400000000001
$ws.Range("M137").Value = -3737.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 39396370
$ws.Range("J2").Value = 4549
$ws.Range("L2").Value = 4549
$ws.Range("N2").Value = -4775

$ws.Range("H5").Value = 35.666668
$ws.Range("I5").Value = 42.4
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 42.4
$ws.Range("L5").Value = 2
$ws.Range("M5").Value = 69.59999999999999
$ws.Range("N5").Value = -226

$ws.Range("H36").Value = 3916.8333
$ws.Range("I36").Value = 3916.8333
$ws.Range("K36").Value = 3916.8333
$ws.Range("M36").Value = -3570.8333

$ws.Range("H45").Value = 3028.4285
$ws.Range("J45").Value = 3199.8333
$ws.Range("L45").Value = 3199.8333
$ws.Range("N45").Value = -3953.8333

$ws.Range("H74").Value = 4149.25
$ws.Range("I74").Value = 3979.6
$ws.Range("K74").Value = 3979.6
$ws.Range("M74").Value = -3105.6

$ws.Range("H77").Value = 4149.25
$ws.Range("I77").Value = 3979.6
$ws.Range("K77").Value = 19898
$ws.Range("M77").Value = -15530

$ws.Range("H102").Value = 18521218
$ws.Range("I102").Value = 55557056
$ws.Range("K102").Value = 55557056
$ws.Range("M102").Value = -55555434

$ws.Range("H110").Value = 3586337.8
$ws.Range("I110").Value = 5557021
$ws.Range("K110").Value = 5557021
$ws.Range("M110").Value = -5554976

$ws.Range("H116").Value = 39396370
$ws.Range("J116").Value = 4549
$ws.Range("L116").Value = 4549
$ws.Range("N116").Value = -9137

$ws.Range("H122").Value = 2541249
$ws.Range("I122").Value = 3369998.8
$ws.Range("J122").Value = 55000
$ws.Range("K122").Value = 10109996.4
$ws.Range("L122").Value = 165000
$ws.Range("M122").Value = -10107546.4
$ws.Range("N122").Value = -169900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 39396370
$ws.Range("J3").Value = 4549
$ws.Range("L3").Value = 4549
$ws.Range("N3").Value = -4777

$ws.Range("H4").Value = 35.666668
$ws.Range("I4").Value = 42.4
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 42.4
$ws.Range("L4").Value = 2
$ws.Range("M4").Value = 72.59999999999999
$ws.Range("N4").Value = -232

$ws.Range("H105").Value = 4171526.8
$ws.Range("I105").Value = 8338064
$ws.Range("K105").Value = 8338064
$ws.Range("M105").Value = -8336317

$ws.Range("H107").Value = 882.2
$ws.Range("I107").Value = 882.2
$ws.Range("K107").Value = 882.2
$ws.Range("M107").Value = 1037.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 150.25
$ws.Range("I5").Value = 150.25
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 150.25
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -38.25
$ws.Range("N5").ClearContents()

$ws.Range("H7").Value = 147.15
$ws.Range("I7").Value = 126.71429
$ws.Range("J7").Value = 194.83333
$ws.Range("K7").Value = 126.71429
$ws.Range("L7").Value = 194.83333
$ws.Range("M7").Value = -13.71429000000001
$ws.Range("N7").Value = -420.83333

$ws.Range("H22").Value = 110174.4
$ws.Range("I22").Value = 111304.78
$ws.Range("K22").Value = 111304.78
$ws.Range("M22").Value = -110954.78

$ws.Range("H25").Value = 4966.6665
$ws.Range("I25").Value = 2450
$ws.Range("K25").Value = 2450
$ws.Range("M25").Value = -2276

$ws.Range("H107").Value = 15152565
$ws.Range("I107").Value = 25000620
$ws.Range("K107").Value = 25000620
$ws.Range("M107").Value = -24998700

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 11111276
$ws.Range("I7").Value = 16666834
$ws.Range("J7").Value = 161.66667
$ws.Range("K7").Value = 50000502
$ws.Range("L7").Value = 485.00001
$ws.Range("M7").Value = -50000390
$ws.Range("N7").Value = -709.00001

$ws.Range("H25").Value = 300
$ws.Range("I25").Value = 300
$ws.Range("K25").Value = 900
$ws.Range("M25").Value = -731

$ws.Range("H30").Value = 300
$ws.Range("I30").Value = 300
$ws.Range("K30").Value = 900
$ws.Range("M30").Value = -798

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 114664.336
$ws.Range("I122").Value = 2832.6667
$ws.Range("K122").Value = 8498.000100000001
$ws.Range("M122").Value = -6048.000100000001

$ws.Range("H132").Value = 1567.3846
$ws.Range("I132").Value = 1567.3846
$ws.Range("K132").Value = 4702.1538
$ws.Range("M132").Value = -2172.1538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 369
$ws.Range("J22").Value = 369
$ws.Range("L22").Value = 369
$ws.Range("N22").Value = -959

$ws.Range("H27").Value = 369
$ws.Range("J27").Value = 369
$ws.Range("L27").Value = 369
$ws.Range("N27").Value = -583

$ws.Range("H30").Value = 2888.3333
$ws.Range("I30").Value = 2082.5
$ws.Range("J30").Value = 4500
$ws.Range("K30").Value = 2082.5
$ws.Range("L30").Value = 4500
$ws.Range("M30").Value = -1974.5
$ws.Range("N30").Value = -4716

$ws.Range("H46").Value = 3272.5454
$ws.Range("I46").Value = 2374.5
$ws.Range("K46").Value = 2374.5
$ws.Range("M46").Value = -2186.5

$ws.Range("H55").Value = 561.8095
$ws.Range("I55").Value = 493.23077
$ws.Range("K55").Value = 493.23077
$ws.Range("M55").Value = -320.23077

$ws.Range("H82").Value = 102199.2
$ws.Range("I82").Value = 2499.125
$ws.Range("J82").Value = 500999.5
$ws.Range("K82").Value = 2499.125
$ws.Range("L82").Value = 500999.5
$ws.Range("M82").Value = -2138.125
$ws.Range("N82").Value = -501721.5

$ws.Range("H85").Value = 102199.2
$ws.Range("I85").Value = 2499.125
$ws.Range("J85").Value = 500999.5
$ws.Range("K85").Value = 2499.125
$ws.Range("L85").Value = 500999.5
$ws.Range("M85").Value = -1251.125
$ws.Range("N85").Value = -503495.5

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4274.273
$ws.Range("I62").Value = 3568
$ws.Range("J62").Value = 4677.857
$ws.Range("K62").Value = 3568
$ws.Range("L62").Value = 4677.857
$ws.Range("M62").Value = -2944
$ws.Range("N62").Value = -5925.857

$ws.Range("H65").Value = 4274.273
$ws.Range("I65").Value = 3568
$ws.Range("J65").Value = 4677.857
$ws.Range("K65").Value = 17840
$ws.Range("L65").Value = 23389.285
$ws.Range("M65").Value = -14720
$ws.Range("N65").Value = -29629.285

$ws.Range("H95").Value = 79999
$ws.Range("J95").Value = 79999
$ws.Range("L95").Value = 79999
$ws.Range("N95").Value = -85491

$ws.Range("H97").Value = 45000
$ws.Range("J97").Value = 45000
$ws.Range("L97").Value = 45000
$ws.Range("N97").Value = -46982

$ws.Range("H122").Value = 3081
$ws.Range("I122").Value = 4900
$ws.Range("J122").Value = 2171.5
$ws.Range("K122").Value = 14700
$ws.Range("L122").Value = 6514.5
$ws.Range("M122").Value = -12250
$ws.Range("N122").Value = -11414.5

$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H126").Value = 1624.25
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = 0

$ws.Range("H132").Value = 3007.182
$ws.Range("I132").Value = 2608.8333
$ws.Range("K132").Value = 7826.499899999999
$ws.Range("M132").Value = -5296.499899999999
